# Add two new trainees (16 certificate rows total) to the DSS Sheet1
# training log, following the exact same per-person block pattern already
# used for the rows immediately above (8 rows per person: one per course).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Copy the existing per-person block formatting (styles only) down ---
# Rows 711-718 use style pair (11 / 45 on the date column); rows 719-726
# use style pair (10 / 44). The two new trainees repeat that same
# alternating pattern, so clone the formats from those template blocks.
$ws.Range("A711:E718").Copy() | Out-Null
$ws.Range("A727:E734").PasteSpecial(-4122) | Out-Null

$ws.Range("A719:E726").Copy() | Out-Null
$ws.Range("A735:E742").PasteSpecial(-4122) | Out-Null

$ws.Application.CutCopyMode = $false

# --- Row height to match the rest of the data rows ---
$ws.Range("A727:E742").RowHeight = 15.75

# --- New shared data ---
$courses = @(
  "30 Hours Construction Safety & Health",
  "30 Hours G. Industry Safety & Health",
  "Electrical Safety ",
  "Fire Marshal",
  "Scaffold Competent Person",
  "Lifting & Rigging Competent Person",
  "Health & Safety Risk Assessment",
  "Safety Management System & PTW"
)
$dates = @(
  "05-12-2024",
  "10-12-2024",
  "06-12-2024",
  "03-12-2024",
  "01-12-2024",
  "02-12-2024",
  "07-12-2024",
  "08-12-2024"
)

$people = @(
  @{ Name = "Seif Elsayed Salem Elsayed"; StartCert = 1726; StartRow = 727 },
  @{ Name = "Ahmed Salah Eldin Ahmed Mohamed"; StartCert = 1734; StartRow = 735 }
)

foreach ($person in $people) {
  for ($i = 0; $i -lt 8; $i++) {
    $r = $person.StartRow + $i
    $certNo = $person.StartCert + $i

    $ws.Range("A$r").Value = "DSS$certNo"
    $ws.Range("B$r").Value = $person.Name
    $ws.Range("C$r").Value = $courses[$i]
    # Leading apostrophe forces text entry so the dd-mm-yyyy string is not
    # auto-converted to a date serial (matches the existing rows, which
    # store the date as plain text too).
    $ws.Range("D$r").Value = "'" + $dates[$i]
    $ws.Range("E$r").Value = 1
  }
}

# --- View state: mirror the scrolled-down selection left by the edit ---
$excel.ActiveWindow.ScrollRow = 725
$ws.Range("G738").Select() | Out-Null

Write-Output "Added rows 727:742"
